$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new literal text value, taken 1:1 from the
# refreshed coinranking.com snapshot (price + 1h volume%, plus the three
# row reorderings at the bottom of the table).
$updates = @(
    @("D2", "26.163.68"),
    @("E2", "  +0.58%  "),
    @("D3", "1.655.43"),
    @("E3", "  +0.44%  "),
    @("E4", "  +0.27%  "),
    @("D5", "217.80"),
    @("E5", "  +0.13%  "),
    @("D6", "0.5293"),
    @("E6", "  +1.53%  "),
    @("D7", "1.006"),
    @("D8", "0.2622"),
    @("E8", "  +0.32%  "),
    @("D9", "0.06318"),
    @("E9", "  +0.76%  "),
    @("D10", "20.39"),
    @("E10", "  -0.34%  "),
    @("D11", "0.07812"),
    @("E11", "  +0.99%  "),
    @("D12", "4.515"),
    @("E12", "  +1.14%  "),
    @("D13", "1.657.72"),
    @("E13", "  +1.01%  "),
    @("D14", "1.883.13"),
    @("E14", "  +0.40%  "),
    @("D15", "0.5486"),
    @("E15", "  +0.92%  "),
    @("D16", "0.0₅8154"),
    @("E16", "  +0.98%  "),
    @("D17", "65.32"),
    @("E17", "  +1.09%  "),
    @("D18", "26.138.03"),
    @("E18", "  +0.43%  "),
    @("E19", "  +0.24%  "),
    @("D20", "4.597"),
    @("E20", "  +0.88%  "),
    @("D21", "190.92"),
    @("E21", "  -0.27%  "),
    @("E22", "  +0.40%  "),
    @("D23", "5.994"),
    @("E23", "  +0.22%  "),
    @("E24", "  +0.25%  "),
    @("D25", "145.48"),
    @("E25", "  +4.80%  "),
    @("D26", "0.1224"),
    @("E26", "  -0.53%  "),
    @("D27", "7.204"),
    @("E27", "  -0.46%  "),
    @("E28", "  -1.08%  "),
    @("D29", "1.472"),
    @("E29", "  +3.21%  "),
    @("D30", "0.05716"),
    @("E30", "  -3.27%  "),
    @("D31", "1.273"),
    @("E31", "  -0.37%  "),
    @("D32", "3.547"),
    @("E32", "  +1.68%  "),
    @("D33", "3.265"),
    @("E33", "  +1.13%  "),
    @("D34", "1.589"),
    @("E34", "  +5.13%  "),
    @("E35", "  +2.02%  "),
    @("D36", "2.421"),
    @("E36", "  +0.24%  "),
    @("D37", "0.9482"),
    @("E37", "  +0.73%  "),
    @("D38", "0.5722"),
    @("E38", "  +0.84%  "),
    @("D39", "0.01607"),
    @("E39", "  +0.22%  "),
    @("D40", "0.8503"),
    @("E40", "  +0.48%  "),
    @("D41", "5.795"),
    @("E41", "  -0.77%  "),
    @("D42", "1.006"),
    @("E42", "  +0.35%  "),
    @("D43", "103.90"),
    @("E43", "  +3.39%  "),
    @("D44", "1.037.36"),
    @("E44", "  +3.72%  "),
    @("D45", "1.795.76"),
    @("E45", "  +0.27%  "),
    @("D46", "56.68"),
    @("B47", "BabyDogeCoin"),
    @("C47", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"),
    @("D47", "0.0₈105"),
    @("E47", "  -2.27%  "),
    @("B48", "Frax"),
    @("C48", "https://coinranking.com/coin/KfWtaeV1W+frax-frax"),
    @("D48", "1.005"),
    @("E48", "  +0.11%  "),
    @("B50", "Cronos"),
    @("C50", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"),
    @("D50", "0.05153"),
    @("E50", "  +0.12%  "),
    @("B51", "EnergySwap"),
    @("C51", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"),
    @("D51", "7.837"),
    @("E51", "  -0.10%  ")
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newValue = $u[1]
    $rng = $ws.Range($cellRef)
    if ($cellRef.StartsWith("D")) {
        # Column D holds price text that can look numeric (e.g. "217.80",
        # "0.05153"); without forcing a Text format first, Excel silently
        # reinterprets it as a floating-point number and mangles the exact
        # string (trailing zeros, precision). Force text, assign, then drop
        # the number format again so no stray style sticks to the cell.
        $rng.NumberFormat = "@"
        $rng.Value = $newValue
        $rng.ClearFormats()
    } else {
        $rng.Value = $newValue
    }
}
